# Update the generated PO test workbook to the new standard template layout.
#
# "Input" sheet: columns are fully reshuffled/renamed (발주일자, 납기일자,
# 거래처명, 거래처 이메일, 납품처명, 납품처 이메일, 프로젝트명, 대분류, 중분류,
# 소분류, 품목명, 규격, 수량, 단가, 총금액, 비고) and shrinks from 17 to 16
# columns (단위/공급가액/부가세/합계 columns are dropped, 비고 becomes empty).
# Header row also loses its bold/bordered style.
#
# "갑지" / "을지" sheets: unchanged data, just drop the leftover empty
# trailing 비고 ("I") cells in rows 2-4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel re-typing it
# as a number/date (e.g. "2025-09-18"), and without leaving a lingering
# explicit cell style behind.
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-NumberCell($ws, $row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

# ---------------------------------------------------------------------------
# Sheet 1: "Input"
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

# Wipe everything (content + formatting) and rebuild from scratch.
$wsInput.Cells.Clear()

$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsInput.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2
Set-TextCell $wsInput 2 1 "2025-09-18"
Set-TextCell $wsInput 2 2 "2025-09-07"
Set-TextCell $wsInput 2 3 "케이에스파워텍"
Set-TextCell $wsInput 2 4 "케이에스파워텍@example.com"
Set-TextCell $wsInput 2 5 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 2 6 "delivery@example.com"
Set-TextCell $wsInput 2 7 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 2 8 "2. 부자재비"
Set-TextCell $wsInput 2 9 "3) 기타"
Set-TextCell $wsInput 2 10 "기타"
Set-TextCell $wsInput 2 11 "스텐망 1083*2145"
Set-TextCell $wsInput 2 12 "KS규격-1"
Set-NumberCell $wsInput 2 13 1
Set-NumberCell $wsInput 2 14 0
Set-NumberCell $wsInput 2 15 0

# Row 3
Set-TextCell $wsInput 3 1 "2025-09-18"
Set-TextCell $wsInput 3 2 "2025-09-08"
Set-TextCell $wsInput 3 3 "케이에스파워텍"
Set-TextCell $wsInput 3 4 "케이에스파워텍@example.com"
Set-TextCell $wsInput 3 5 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 3 6 "delivery@example.com"
Set-TextCell $wsInput 3 7 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 3 8 "2. 부자재비"
Set-TextCell $wsInput 3 9 "3) 기타"
Set-TextCell $wsInput 3 10 "기타"
Set-TextCell $wsInput 3 11 "스텐망 1168*343"
Set-TextCell $wsInput 3 12 "KS규격-2"
Set-NumberCell $wsInput 3 13 8
Set-NumberCell $wsInput 3 14 29000
Set-NumberCell $wsInput 3 15 255200

# Row 4
Set-TextCell $wsInput 4 1 "2025-09-11"
Set-TextCell $wsInput 4 2 "2025-09-19"
Set-TextCell $wsInput 4 3 "케이에스파워텍"
Set-TextCell $wsInput 4 4 "케이에스파워텍@example.com"
Set-TextCell $wsInput 4 5 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 4 6 "delivery@example.com"
Set-TextCell $wsInput 4 7 "힐스테이트 도곡동1차"
Set-TextCell $wsInput 4 8 "5. 운반비"
Set-TextCell $wsInput 4 9 "일반자재"
Set-TextCell $wsInput 4 10 "기타"
Set-TextCell $wsInput 4 11 "3월 운반비"
Set-TextCell $wsInput 4 12 "KS규격-3"
Set-NumberCell $wsInput 4 13 1
Set-NumberCell $wsInput 4 14 0
Set-NumberCell $wsInput 4 15 0

# ---------------------------------------------------------------------------
# Sheets 2 & 3: "갑지" / "을지" - just drop the leftover empty 비고 cells
# in column I (rows 2-4).
# ---------------------------------------------------------------------------
foreach ($sheetName in @("갑지", "을지")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("I2").ClearContents()
    $ws.Range("I3").ClearContents()
    $ws.Range("I4").ClearContents()
}
